# Insert a new data row at row 114 (this pushes the existing rows 114..194
# down to 115..195, preserving all of their data/styles), then populate the
# newly inserted row 114 with the new "Arveja Verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 114, shifting rows 114-194 down
# to 115-195.
$ws.Rows(114).Insert()

# Populate the new row 114 with the new record's data.
$ws.Cells.Item(114, 1).Value  = 9
$ws.Cells.Item(114, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(114, 3).Value  = "Metropolitana"
$ws.Cells.Item(114, 4).Value  = 45233
$ws.Cells.Item(114, 5).Value  = 13
$ws.Cells.Item(114, 6).Value  = 100112022
$ws.Cells.Item(114, 7).Value  = "Arveja Verde"
$ws.Cells.Item(114, 8).Value  = "Sin especificar"
$ws.Cells.Item(114, 9).Value  = "Primera"
$ws.Cells.Item(114, 10).Value = 52
$ws.Cells.Item(114, 11).Value = 17000
$ws.Cells.Item(114, 12).Value = 19000
$ws.Cells.Item(114, 13).Value = 18000
$ws.Cells.Item(114, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(114, 16).Value = 720
$ws.Cells.Item(114, 17).Value = 25
$ws.Cells.Item(114, 18).Value = "Hortaliza"
